# Applies the "INS final commit remaining unstaged changes" edits to
# TC10_INS_CancerType-Leukemia.xlsx:
#  - rewrites the 4 SQL statements stored in column B (Programs / Projects /
#    Grants / Publications query text) to the corrected versions
#  - normalizes the formatting of B2:B5 and C2 so they all share one
#    consistent font/wrap style (mirrors the font/style table cleanup in the
#    diff, where several duplicate font entries collapsed into one)
#  - resets the sheet view: no more scrolled topLeftCell, and the active
#    selection moves from C8 to C2

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Updated SQL text for the four query cells
# ---------------------------------------------------------------------

$programsQuery = @'
SELECT DISTINCT 
    prg.program_name AS "Program",
  CASE
    WHEN prg.program_link IS NOT NULL THEN prg.program_acronym
        ELSE prg.program_link
    END  AS "Website",
    prg.focus_area AS "Special Topic",
    prg.cancer_type AS "Cancer Type",
 CASE 
        WHEN prg.data_link IS NOT NULL THEN prg.program_acronym     
        ELSE prg.data_link
    END AS "Data Location Details" 
FROM 
    df_program prg
WHERE 
     prg.cancer_type LIKE '%Leukemia%'
ORDER BY 
    lower(prg.program_name) ASC
LIMIT 100;
'@

$projectsQuery = @'
SELECT DISTINCT
    prj.project_id AS "Project ID", 
    prj.project_title AS "Project Title",
    prj.project_org_name AS "Organization",
    prj.project_start_date AS "Project Start Date",
    prj.project_end_date AS "Project End Date"
FROM 
    df_project prj
LEFT JOIN 
    df_program prg ON prj."program.program_id" = prg.program_id
LEFT JOIN 
    df_grant gnt ON prj.project_id = gnt."project.project_id"
LEFT JOIN 
    df_publication pub ON prj.project_id = pub."project.project_id"
WHERE 
     prg.cancer_type LIKE '%Leukemia%'
ORDER BY 
    lower(prj.project_id) ASC
LIMIT 100;
'@

$grantsQuery = @'
SELECT DISTINCT
    gnt.grant_id AS "Grant ID", 
    prj.project_id AS "Project",
    gnt.grant_title AS "Grant Title",
    gnt.principal_investigators AS "Principal Investigators",
    gnt.program_officers AS "Program Officers",
    gnt.fiscal_year AS "Fiscal Year",
    gnt.grant_end_date AS "Project End Date"
FROM 
    df_grant gnt
LEFT JOIN 
    df_project prj ON gnt."project.project_id" = prj.project_id
LEFT JOIN 
    df_program prg ON prj."program.program_id" = prg.program_id
LEFT JOIN 
    df_publication pub ON prj.project_id = pub."project.project_id"
WHERE 
    prg.cancer_type  LIKE '%Leukemia%'
ORDER BY 
    lower(gnt.grant_id) ASC
LIMIT 100;
'@

$publicationsQuery = @'
SELECT DISTINCT
    pub.pmid AS "PubMed ID", 
    pub.publication_title AS "Title",
    pub.authors AS "Authors",
    pub.publication_date AS "Publication Date",
    pub.cited_by AS "Cited By",
    CASE 
    WHEN pub.relative_citation_ratio = 0 THEN '0'
    WHEN pub.relative_citation_ratio = 7.0 THEN '7'
    WHEN pub.relative_citation_ratio = 2.0 THEN '2'
  WHEN pub.relative_citation_ratio = 1.0 THEN '1'
    WHEN pub.relative_citation_ratio = ROUND(pub.relative_citation_ratio) THEN CAST(ROUND(pub.relative_citation_ratio) AS VARCHAR) 
    ELSE CAST(ROUND(pub.relative_citation_ratio, 2) AS VARCHAR)
END AS "Relative Citation Ratio"
FROM 
    df_publication pub
LEFT JOIN 
    df_project prj ON pub."project.project_id" = prj.project_id
LEFT JOIN 
    df_program prg ON prj."program.program_id" = prg.program_id
LEFT JOIN 
    df_grant gnt ON prj.project_id = gnt."project.project_id"
WHERE 
     prg.cancer_type  LIKE '%Leukemia%'
ORDER BY 
    lower(pub.pmid) ASC
LIMIT 100;
'@


$ws.Range("B2").Value = $programsQuery
$ws.Range("B3").Value = $projectsQuery
$ws.Range("B4").Value = $grantsQuery
$ws.Range("B5").Value = $publicationsQuery

# ---------------------------------------------------------------------
# 2. Normalize formatting on B2:B5 and C2 to one shared style
#    (matches the font/cellXf de-duplication seen in styles.xml)
# ---------------------------------------------------------------------

$bRange = $ws.Range("B2:B5")
$bRange.WrapText = $true
$bRange.Font.Size = 12

$cRange = $ws.Range("C2")
$cRange.WrapText = $true
$cRange.Font.Size = 12

# ---------------------------------------------------------------------
# 3. Reset the view: select C2 (was C8), dropping the scrolled
#    topLeftCell="A5" state from the sheet view
# ---------------------------------------------------------------------

$ws.Range("C2").Select()
